$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 198, shifting old rows 198-293 down to 200-295
$ws.Range("A198:A199").EntireRow.Insert()

# Set values for new row 198
$ws.Cells.Item(198,1).Value = 10
$ws.Cells.Item(198,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(198,3).Value = "La Araucanía"
$ws.Cells.Item(198,4).Value = 44510
$ws.Cells.Item(198,5).Value = 9
$ws.Cells.Item(198,6).Value = 100112023
$ws.Cells.Item(198,7).Value = "Brócoli"
$ws.Cells.Item(198,8).Value = "Sin especificar"
$ws.Cells.Item(198,9).Value = "Primera"
$ws.Cells.Item(198,10).Value = 1400
$ws.Cells.Item(198,11).Value = 800
$ws.Cells.Item(198,12).Value = 900
$ws.Cells.Item(198,13).Value = 843
$ws.Cells.Item(198,14).Value = "`$/unidad"
$ws.Cells.Item(198,15).Value = "Región Metropolitana"
$ws.Cells.Item(198,16).Value = 843
$ws.Cells.Item(198,17).Value = 1
$ws.Cells.Item(198,18).Value = "Hortaliza"

# Set values for new row 199
$ws.Cells.Item(199,1).Value = 10
$ws.Cells.Item(199,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(199,3).Value = "La Araucanía"
$ws.Cells.Item(199,4).Value = 44510
$ws.Cells.Item(199,5).Value = 9
$ws.Cells.Item(199,6).Value = 100112023
$ws.Cells.Item(199,7).Value = "Brócoli"
$ws.Cells.Item(199,8).Value = "Sin especificar"
$ws.Cells.Item(199,9).Value = "Primera"
$ws.Cells.Item(199,10).Value = 800
$ws.Cells.Item(199,11).Value = 900
$ws.Cells.Item(199,12).Value = 900
$ws.Cells.Item(199,13).Value = 900
$ws.Cells.Item(199,14).Value = "`$/unidad"
$ws.Cells.Item(199,15).Value = "Región de O'Higgins"
$ws.Cells.Item(199,16).Value = 900
$ws.Cells.Item(199,17).Value = 1
$ws.Cells.Item(199,18).Value = "Hortaliza"
